{"js": "const pairs = [\n  [\"90\u00f73=30, 0\", \"24\u00f77=3, 3\"],\n  [\"85\u00f76=14, 1\", \"82\u00f78=10, 2\"],\n  [\"96\u00f77=13, 5\", \"27\u00f78=3, 3\"],\n  [\"22\u00f73=7, 1\", \"79\u00f75=15, 4\"],\n  [\"85\u00f72=42, 1\", \"54\u00f77=7, 5\"],\n  [\"21\u00f72=10, 1\", \"23\u00f76=3, 5\"],\n  [\"36\u00f74=9, 0\", \"45\u00f73=15, 0\"],\n  [\"88\u00f75=17, 3\", \"43\u00f76=7, 1\"],\n  [\"65\u00f73=21, 2\", \"43\u00f74=10, 3\"],\n  [\"29\u00f74=7, 1\", \"25\u00f76=4, 1\"],\n  [\"80\u00f74=20, 0\", \"51\u00f73=17, 0\"],\n  [\"16\u00f72=8, 0\", \"11\u00f74=2, 3\"],\n  [\"30\u00f76=5, 0\", \"10\u00f72=5, 0\"],\n  [\"48\u00f72=24, 0\", \"68\u00f72=34, 0\"],\n  [\"82\u00f73=27, 1\", \"52\u00f76=8, 4\"],\n  [\"57\u00f74=14, 1\", \"52\u00f75=10, 2\"],\n  [\"19\u00f78=2, 3\", \"49\u00f76=8, 1\"],\n  [\"71\u00f73=23, 2\", \"88\u00f75=17, 3\"],\n  [\"71\u00f77=10, 1\", \"20\u00f74=5, 0\"],\n  [\"93\u00f75=18, 3\", \"36\u00f75=7, 1\"],\n  [\"87\u00f79=9, 6\", \"24\u00f74=6, 0\"],\n  [\"36\u00f72=18, 0\", \"75\u00f79=8, 3\"],\n  [\"64\u00f79=7, 1\", \"28\u00f72=14, 0\"],\n  [\"66\u00f76=11, 0\", \"24\u00f73=8, 0\"],\n  [\"59\u00f73=19, 2\", \"95\u00f77=13, 4\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"90\u00f73=30, 0\", \"24\u00f77=3, 3\"),\n    @(\"85\u00f76=14, 1\", \"82\u00f78=10, 2\"),\n    @(\"96\u00f77=13, 5\", \"27\u00f78=3, 3\"),\n    @(\"22\u00f73=7, 1\", \"79\u00f75=15, 4\"),\n    @(\"85\u00f72=42, 1\", \"54\u00f77=7, 5\"),\n    @(\"21\u00f72=10, 1\", \"23\u00f76=3, 5\"),\n    @(\"36\u00f74=9, 0\", \"45\u00f73=15, 0\"),\n    @(\"88\u00f75=17, 3\", \"43\u00f76=7, 1\"),\n    @(\"65\u00f73=21, 2\", \"43\u00f74=10, 3\"),\n    @(\"29\u00f74=7, 1\", \"25\u00f76=4, 1\"),\n    @(\"80\u00f74=20, 0\", \"51\u00f73=17, 0\"),\n    @(\"16\u00f72=8, 0\", \"11\u00f74=2, 3\"),\n    @(\"30\u00f76=5, 0\", \"10\u00f72=5, 0\"),\n    @(\"48\u00f72=24, 0\", \"68\u00f72=34, 0\"),\n    @(\"82\u00f73=27, 1\", \"52\u00f76=8, 4\"),\n    @(\"57\u00f74=14, 1\", \"52\u00f75=10, 2\"),\n    @(\"19\u00f78=2, 3\", \"49\u00f76=8, 1\"),\n    @(\"71\u00f73=23, 2\", \"88\u00f75=17, 3\"),\n    @(\"71\u00f77=10, 1\", \"20\u00f74=5, 0\"),\n    @(\"93\u00f75=18, 3\", \"36\u00f75=7, 1\"),\n    @(\"87\u00f79=9, 6\", \"24\u00f74=6, 0\"),\n    @(\"36\u00f72=18, 0\", \"75\u00f79=8, 3\"),\n    @(\"64\u00f79=7, 1\", \"28\u00f72=14, 0\"),\n    @(\"66\u00f76=11, 0\", \"24\u00f73=8, 0\"),\n    @(\"59\u00f73=19, 2\", \"95\u00f77=13, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
